# Update gh-pages output (合肥 漫展信息) to the state generated at 456a3b4.
#
# Touches the two sheets that hold the 展览 ("exhibitions") table data —
# "展览" and "全部类型" — which were kept in sync in the source data:
#   * bump a handful of "want to go" head-counts (column F / some column G)
#   * insert a brand-new event row (2024-05-03, 合肥·百合Only2.0·同好交流)
#     before the existing "2024-05-18" row, shifting the two rows below it
#     down by one
#   * bump the "2024-06-01" row's go-count now that it has moved down

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- 1. simple numeric refreshes on existing rows ------------------
    $ws.Range("F3").Value = 7271
    $ws.Range("F4").Value = 5462
    $ws.Range("F6").Value = 169
    $ws.Range("F8").Value = 43
    $ws.Range("F13").Value = 22
    $ws.Range("F15").Value = 258

    if ($sheetName -eq "展览") {
        $ws.Range("F11").Value = 101
    } else {
        $ws.Range("F11").Value = 102
    }

    # --- 2. insert the new row 17 (shifts old 17->18, old 18->19) ------
    $ws.Rows(17).Insert()

    # Row 17 lands empty after the insert; seed its look from the row
    # that is now directly below it (the old row 17) so borders/fonts on
    # column A etc. match the rest of the table.
    $ws.Range("A18:I18").Copy()
    $ws.Range("A17:I17").PasteSpecial(-4122)

    # B17 must hold literal text ("2024-05-03"), not an auto-converted
    # date serial -- force text format before assigning it, then restore
    # the plain (unformatted) look used by every other date cell in the
    # column by re-pasting that cell's format from its neighbour.
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("A17").Value = 16
    $ws.Range("B17").Value = "2024-05-03"
    $ws.Range("C17").Value = "合肥·百合Only2.0·同好交流"
    $ws.Range("D17").Value = "北二环与新蚌埠路交汇处 蓝金湾大酒店"
    $ws.Range("E17").Value = "2024.05.03 10:00-05.03 16:00"
    $ws.Range("F17").Value = 5
    $ws.Range("G17").Value = 46
    $ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=83045"
    $ws.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202403/wDDNHQGa1710906388610.jpeg"

    $ws.Range("B18").Copy()
    $ws.Range("B17").PasteSpecial(-4122)

    # --- 3. the two rows pushed down by the insert keep their data, but
    #        their sequential index in column A (= row-1) has to be
    #        renumbered, and the old row 18 (now row 19) gets a
    #        refreshed go-count ----------------------------------------
    $ws.Range("A18").Value = 17
    $ws.Range("A19").Value = 18
    $ws.Range("F19").Value = 29
}
